# Auto-generated edit script: updates market-price-derived columns (H-N)
# across the Leve profit tables in each job sheet, per scheduled price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 416855.16
$ws.Range("J9").Value = 714394.5600000001
$ws.Range("L9").Value = 714394.5600000001
$ws.Range("N9").Value = -714732.5600000001
$ws.Range("H42").Value = 297.45456
$ws.Range("I42").Value = 156.71428
$ws.Range("J42").Value = 543.75
$ws.Range("K42").Value = 470.14284
$ws.Range("L42").Value = 1631.25
$ws.Range("M42").Value = -240.14284
$ws.Range("N42").Value = -2091.25
$ws.Range("H55").Value = 129.14285
$ws.Range("J55").Value = 113.3
$ws.Range("L55").Value = 113.3
$ws.Range("N55").Value = -541.3
$ws.Range("H111").Value = 2327.5715
$ws.Range("I111").Value = 2382.1667
$ws.Range("K111").Value = 7146.500100000001
$ws.Range("M111").Value = -4079.500100000001
$ws.Range("H132").Value = 3298.5
$ws.Range("I132").Value = 3298.5
$ws.Range("K132").Value = 9895.5
$ws.Range("M132").Value = -7365.5
$ws.Range("H133").Value = 67500
$ws.Range("J133").Value = 67500
$ws.Range("L133").Value = 67500
$ws.Range("N133").Value = -77620
$ws.Range("H137").Value = 2497.36
$ws.Range("I137").Value = 1429.7333
$ws.Range("J137").Value = 4098.8
$ws.Range("K137").Value = 4289.199900000001
$ws.Range("L137").Value = 12296.4
$ws.Range("M137").Value = -1739.199900000001
$ws.Range("N137").Value = -17396.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4448.7144
$ws.Range("I32").Value = 2350.0571
$ws.Range("J32").Value = 14942
$ws.Range("K32").Value = 2350.0571
$ws.Range("L32").Value = 14942
$ws.Range("M32").Value = -2063.0571
$ws.Range("N32").Value = -15516
$ws.Range("H61").Value = 4505.913
$ws.Range("I61").Value = 4685.9443
$ws.Range("K61").Value = 4685.9443
$ws.Range("M61").Value = -4473.9443
$ws.Range("H63").Value = 2747.2856
$ws.Range("I63").Value = 2532.9
$ws.Range("J63").Value = 3283.25
$ws.Range("K63").Value = 2532.9
$ws.Range("L63").Value = 3283.25
$ws.Range("M63").Value = -1846.9
$ws.Range("N63").Value = -4655.25
$ws.Range("H66").Value = 2747.2856
$ws.Range("I66").Value = 2532.9
$ws.Range("J66").Value = 3283.25
$ws.Range("K66").Value = 12664.5
$ws.Range("L66").Value = 16416.25
$ws.Range("M66").Value = -9232.5
$ws.Range("N66").Value = -23280.25
$ws.Range("H122").Value = 2213
$ws.Range("I122").Value = 2255.4285
$ws.Range("K122").Value = 6766.2855
$ws.Range("M122").Value = -4316.2855
$ws.Range("H136").Value = 4505.913
$ws.Range("I136").Value = 4685.9443
$ws.Range("K136").Value = 14057.8329
$ws.Range("M136").Value = -11507.8329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2453.2778
$ws.Range("I86").Value = 2424
$ws.Range("J86").Value = 2499.2856
$ws.Range("K86").Value = 2424
$ws.Range("L86").Value = 2499.2856
$ws.Range("M86").Value = -1301
$ws.Range("N86").Value = -4745.2856
$ws.Range("H89").Value = 2453.2778
$ws.Range("I89").Value = 2424
$ws.Range("J89").Value = 2499.2856
$ws.Range("K89").Value = 12120
$ws.Range("L89").Value = 12496.428
$ws.Range("M89").Value = -6504
$ws.Range("N89").Value = -23728.428
$ws.Range("H94").Value = 829.2727
$ws.Range("I94").Value = 723.4737
$ws.Range("J94").Value = 1499.3334
$ws.Range("K94").Value = 723.4737
$ws.Range("L94").Value = 1499.3334
$ws.Range("M94").Value = -272.4737
$ws.Range("N94").Value = -2401.3334
$ws.Range("H105").Value = 2359.111
$ws.Range("I105").Value = 2615.8125
$ws.Range("J105").Value = 305.5
$ws.Range("K105").Value = 2615.8125
$ws.Range("L105").Value = 305.5
$ws.Range("M105").Value = -868.8125
$ws.Range("N105").Value = -3799.5
$ws.Range("H107").Value = 2115.6667
$ws.Range("I107").Value = 1687.6666
$ws.Range("K107").Value = 1687.6666
$ws.Range("M107").Value = 232.3334
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("H134").Value = 3365.1428
$ws.Range("I134").Value = 3016.7144
$ws.Range("K134").Value = 9050.143199999999
$ws.Range("M134").Value = -6515.143199999999
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4226.625
$ws.Range("I31").Value = 3755.4167
$ws.Range("K31").Value = 3755.4167
$ws.Range("M31").Value = -3460.4167
$ws.Range("H34").Value = 4226.625
$ws.Range("I34").Value = 3755.4167
$ws.Range("K34").Value = 3755.4167
$ws.Range("M34").Value = -3553.4167

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10765288
$ws.Range("I4").Value = 2010379.8
$ws.Range("K4").Value = 6031139.4
$ws.Range("M4").Value = -6031027.4
$ws.Range("H22").Value = 775
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 2400
$ws.Range("L22").Value = 2250
$ws.Range("M22").Value = -2231
$ws.Range("N22").Value = -2588
$ws.Range("H27").Value = 775
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 750
$ws.Range("K27").Value = 2400
$ws.Range("L27").Value = 2250
$ws.Range("M27").Value = -2298
$ws.Range("N27").Value = -2454
$ws.Range("H51").Value = 2158.1667
$ws.Range("I51").Value = 2146.25
$ws.Range("J51").Value = 2182
$ws.Range("K51").Value = 6438.75
$ws.Range("L51").Value = 6546
$ws.Range("M51").Value = -5978.75
$ws.Range("N51").Value = -7466
$ws.Range("H75").Value = 815
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 815
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 2445
$ws.Range("N75").Value = -4441
$ws.Range("H78").Value = 815
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 815
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 7335
$ws.Range("N78").Value = -17319
$ws.Range("H92").Value = 289.25
$ws.Range("J92").Value = 289.33334
$ws.Range("L92").Value = 868.0000200000001
$ws.Range("N92").Value = -3364.00002
$ws.Range("H112").Value = 480.33334
$ws.Range("I112").Value = 480.33334
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 1441.00002
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -333.0000199999999
$ws.Range("H114").Value = 1285.25
$ws.Range("I114").Value = 1285.25
$ws.Range("K114").Value = 3855.75
$ws.Range("M114").Value = -601.75
$ws.Range("H117").Value = 806.2
$ws.Range("I117").Value = 899.75
$ws.Range("J117").Value = 432
$ws.Range("K117").Value = 2699.25
$ws.Range("L117").Value = 1296
$ws.Range("M117").Value = 742.75
$ws.Range("N117").Value = -8180
$ws.Range("H129").Value = 4386.2
$ws.Range("I129").Value = 976.5
$ws.Range("J129").Value = 5626.091
$ws.Range("K129").Value = 2929.5
$ws.Range("L129").Value = 16878.273
$ws.Range("M129").Value = 2070.5
$ws.Range("N129").Value = -26878.273
$ws.Range("M75").ClearContents()
$ws.Range("M78").ClearContents()
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 157.6842
$ws.Range("I2").Value = 177.26666
$ws.Range("J2").Value = 84.25
$ws.Range("K2").Value = 177.26666
$ws.Range("L2").Value = 84.25
$ws.Range("M2").Value = -64.26666
$ws.Range("N2").Value = -310.25
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("H122").Value = 2929.9412
$ws.Range("I122").Value = 3189.75
$ws.Range("K122").Value = 9569.25
$ws.Range("M122").Value = -7119.25
$ws.Range("M18").ClearContents()
$ws.Range("N51").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1234.8572
$ws.Range("I22").Value = 1786.75
$ws.Range("K22").Value = 1786.75
$ws.Range("M22").Value = -1491.75
$ws.Range("H27").Value = 1234.8572
$ws.Range("I27").Value = 1786.75
$ws.Range("K27").Value = 1786.75
$ws.Range("M27").Value = -1679.75
$ws.Range("H40").Value = 5105.7095
$ws.Range("I40").Value = 4697.35
$ws.Range("K40").Value = 4697.35
$ws.Range("M40").Value = -4561.35
$ws.Range("H119").Value = 102367.8
$ws.Range("J119").Value = 102367.8
$ws.Range("L119").Value = 102367.8
$ws.Range("N119").Value = -112043.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
$ws.Range("H132").Value = 5325
$ws.Range("I132").Value = 5325
$ws.Range("K132").Value = 15975
$ws.Range("M132").Value = -13445
